# Apply edits to "Test 2.xlsx":
#  - Header row 3: D3 "BS" -> "BD", add E3 "children"
#  - Data row 4: C4 numeric 4 -> text "5", add E4 = 1
#  - Data row 5: C5 numeric 1.8 -> text "10", add E5 = 2
#  - New rows 6-9 with only column E populated (2, 3, 3, 3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (second header row) ---
# Set new string values first (order chosen so shared strings are appended as
# "10", "5", "children" before "BD", matching the target layout as closely as possible).
$ws.Range("C5").Value = "10"
$ws.Range("C4").Value = "5"
$ws.Range("E3").Value = "children"
$ws.Range("D3").Value = "BD"

# --- Row 4 / Row 5 numeric additions ---
# Set numeric values BEFORE copying the number format so that they remain
# true numeric cells (no t="s") even though the column uses a text ("@") format.
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 2

# --- New rows 6-9, column E only ---
$ws.Range("E6").Value = 2
$ws.Range("E7").Value = 3
$ws.Range("E8").Value = 3
$ws.Range("E9").Value = 3

# Apply the same style/number-format used by the rest of the table (style index 1,
# text format "@") to the newly touched cells, matching existing column styling.
$refFormat = $ws.Range("A5").NumberFormat
$ws.Range("E3:E9").NumberFormat = $refFormat
$ws.Range("A6:D9").NumberFormat = $refFormat

Write-Host "Edit applied."
